# ===================================================================
# England Championship 2023-2024 - apply committed edit
# Part 1: reorder match rows 15-23 and 66-76 (data re-sequenced by date)
# ===================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: was match #14 content, now gets match data previously stored at row 22
$ws.Cells.Item(15, 6).Value = "Cardiff"
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = "QPR"
$ws.Cells.Item(15, 9).Value = 2
$ws.Cells.Item(15, 10).Value = 2.18
$ws.Cells.Item(15, 11).Value = "06/08/2023 15:42"
$ws.Cells.Item(15, 12).Value = 1.91
$ws.Cells.Item(15, 13).Value = "12/08/2023 15:41"
$ws.Cells.Item(15, 14).Value = 3.25
$ws.Cells.Item(15, 15).Value = "06/08/2023 15:42"
$ws.Cells.Item(15, 16).Value = 3.64
$ws.Cells.Item(15, 17).Value = "12/08/2023 15:57"
$ws.Cells.Item(15, 18).Value = 3.78
$ws.Cells.Item(15, 19).Value = "06/08/2023 15:42"
$ws.Cells.Item(15, 20).Value = 4.32
$ws.Cells.Item(15, 21).Value = "12/08/2023 15:57"
$ws.Cells.Item(15, 22).Value = "https://www.betexplorer.com/football/england/championship/cardiff-qpr/8tI21653/"

# Row 16: was match #15 content, now gets match data previously stored at row 20
$ws.Cells.Item(16, 6).Value = "Huddersfield"
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = "Leicester"
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 10).Value = 5.17
$ws.Cells.Item(16, 11).Value = "07/08/2023 13:07"
$ws.Cells.Item(16, 12).Value = 4.02
$ws.Cells.Item(16, 13).Value = "12/08/2023 15:59"
$ws.Cells.Item(16, 14).Value = 4.25
$ws.Cells.Item(16, 15).Value = "07/08/2023 13:07"
$ws.Cells.Item(16, 16).Value = 3.76
$ws.Cells.Item(16, 17).Value = "12/08/2023 15:55"
$ws.Cells.Item(16, 18).Value = 1.65
$ws.Cells.Item(16, 19).Value = "07/08/2023 13:07"
$ws.Cells.Item(16, 20).Value = 1.94
$ws.Cells.Item(16, 21).Value = "12/08/2023 15:55"
$ws.Cells.Item(16, 22).Value = "https://www.betexplorer.com/football/england/championship/huddersfield-leicester/2kUPRg0P/"

# Row 17: was match #16 content, now gets match data previously stored at row 19
$ws.Cells.Item(17, 6).Value = "Hull"
$ws.Cells.Item(17, 7).Value = 4
$ws.Cells.Item(17, 8).Value = "Sheffield Wed"
$ws.Cells.Item(17, 9).Value = 2
$ws.Cells.Item(17, 10).Value = 2.4
$ws.Cells.Item(17, 11).Value = "07/08/2023 13:08"
$ws.Cells.Item(17, 12).Value = 2.47
$ws.Cells.Item(17, 13).Value = "12/08/2023 15:54"
$ws.Cells.Item(17, 14).Value = 3.46
$ws.Cells.Item(17, 15).Value = "07/08/2023 13:08"
$ws.Cells.Item(17, 16).Value = 3.25
$ws.Cells.Item(17, 17).Value = "12/08/2023 15:50"
$ws.Cells.Item(17, 18).Value = 3.05
$ws.Cells.Item(17, 19).Value = "07/08/2023 13:08"
$ws.Cells.Item(17, 20).Value = 3.18
$ws.Cells.Item(17, 21).Value = "12/08/2023 15:54"
$ws.Cells.Item(17, 22).Value = "https://www.betexplorer.com/football/england/championship/hull-city-sheffield-wed/4bGiND8t/"

# Row 19: was match #18 content, now gets match data previously stored at row 16
$ws.Cells.Item(19, 6).Value = "Preston"
$ws.Cells.Item(19, 7).Value = 2
$ws.Cells.Item(19, 8).Value = "Sunderland"
$ws.Cells.Item(19, 9).Value = 1
$ws.Cells.Item(19, 10).Value = 3.02
$ws.Cells.Item(19, 11).Value = "06/08/2023 18:13"
$ws.Cells.Item(19, 12).Value = 2.54
$ws.Cells.Item(19, 13).Value = "12/08/2023 15:52"
$ws.Cells.Item(19, 14).Value = 3.42
$ws.Cells.Item(19, 15).Value = "06/08/2023 18:13"
$ws.Cells.Item(19, 16).Value = 3.36
$ws.Cells.Item(19, 17).Value = "12/08/2023 15:29"
$ws.Cells.Item(19, 18).Value = 2.44
$ws.Cells.Item(19, 19).Value = "06/08/2023 18:13"
$ws.Cells.Item(19, 20).Value = 2.97
$ws.Cells.Item(19, 21).Value = "12/08/2023 15:52"
$ws.Cells.Item(19, 22).Value = "https://www.betexplorer.com/football/england/championship/preston-sunderland/fm94KBha/"

# Row 20: was match #19 content, now gets match data previously stored at row 15
$ws.Cells.Item(20, 6).Value = "Rotherham"
$ws.Cells.Item(20, 7).Value = 2
$ws.Cells.Item(20, 8).Value = "Blackburn"
$ws.Cells.Item(20, 9).Value = 2
$ws.Cells.Item(20, 10).Value = 2.86
$ws.Cells.Item(20, 11).Value = "05/08/2023 16:12"
$ws.Cells.Item(20, 12).Value = 3.48
$ws.Cells.Item(20, 13).Value = "12/08/2023 15:57"
$ws.Cells.Item(20, 14).Value = 3.26
$ws.Cells.Item(20, 15).Value = "05/08/2023 16:12"
$ws.Cells.Item(20, 16).Value = 3.48
$ws.Cells.Item(20, 17).Value = "12/08/2023 15:54"
$ws.Cells.Item(20, 18).Value = 2.68
$ws.Cells.Item(20, 19).Value = "05/08/2023 16:12"
$ws.Cells.Item(20, 20).Value = 2.2
$ws.Cells.Item(20, 21).Value = "12/08/2023 15:57"
$ws.Cells.Item(20, 22).Value = "https://www.betexplorer.com/football/england/championship/rotherham-blackburn/Yw88JV75/"

# Row 21: was match #20 content, now gets match data previously stored at row 23
$ws.Cells.Item(21, 6).Value = "Southampton"
$ws.Cells.Item(21, 7).Value = 4
$ws.Cells.Item(21, 8).Value = "Norwich"
$ws.Cells.Item(21, 9).Value = 4
$ws.Cells.Item(21, 10).Value = 1.79
$ws.Cells.Item(21, 11).Value = "07/08/2023 13:08"
$ws.Cells.Item(21, 12).Value = 2.04
$ws.Cells.Item(21, 13).Value = "12/08/2023 15:55"
$ws.Cells.Item(21, 14).Value = 3.82
$ws.Cells.Item(21, 15).Value = "07/08/2023 13:08"
$ws.Cells.Item(21, 16).Value = 3.69
$ws.Cells.Item(21, 17).Value = "12/08/2023 15:55"
$ws.Cells.Item(21, 18).Value = 4.66
$ws.Cells.Item(21, 19).Value = "07/08/2023 13:08"
$ws.Cells.Item(21, 20).Value = 3.74
$ws.Cells.Item(21, 21).Value = "12/08/2023 15:55"
$ws.Cells.Item(21, 22).Value = "https://www.betexplorer.com/football/england/championship/southampton-norwich/jF2DIkNB/"

# Row 22: was match #21 content, now gets match data previously stored at row 21
$ws.Cells.Item(22, 6).Value = "Birmingham"
$ws.Cells.Item(22, 7).Value = 1
$ws.Cells.Item(22, 8).Value = "Leeds"
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 4.3
$ws.Cells.Item(22, 11).Value = "07/08/2023 13:07"
$ws.Cells.Item(22, 12).Value = 3.37
$ws.Cells.Item(22, 13).Value = "12/08/2023 15:58"
$ws.Cells.Item(22, 14).Value = 3.84
$ws.Cells.Item(22, 15).Value = "07/08/2023 13:07"
$ws.Cells.Item(22, 16).Value = 3.3
$ws.Cells.Item(22, 17).Value = "12/08/2023 15:48"
$ws.Cells.Item(22, 18).Value = 1.84
$ws.Cells.Item(22, 19).Value = "07/08/2023 13:07"
$ws.Cells.Item(22, 20).Value = 2.34
$ws.Cells.Item(22, 21).Value = "12/08/2023 15:58"
$ws.Cells.Item(22, 22).Value = "https://www.betexplorer.com/football/england/championship/birmingham-leeds/dhsv4j8I/"

# Row 23: was match #22 content, now gets match data previously stored at row 17
$ws.Cells.Item(23, 6).Value = "Millwall"
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = "Bristol City"
$ws.Cells.Item(23, 9).Value = 1
$ws.Cells.Item(23, 10).Value = 1.71
$ws.Cells.Item(23, 11).Value = "05/08/2023 16:12"
$ws.Cells.Item(23, 12).Value = 1.93
$ws.Cells.Item(23, 13).Value = "12/08/2023 15:59"
$ws.Cells.Item(23, 14).Value = 3.82
$ws.Cells.Item(23, 15).Value = "05/08/2023 16:12"
$ws.Cells.Item(23, 16).Value = 3.57
$ws.Cells.Item(23, 17).Value = "12/08/2023 15:59"
$ws.Cells.Item(23, 18).Value = 5.37
$ws.Cells.Item(23, 19).Value = "05/08/2023 16:12"
$ws.Cells.Item(23, 20).Value = 4.34
$ws.Cells.Item(23, 21).Value = "12/08/2023 15:59"
$ws.Cells.Item(23, 22).Value = "https://www.betexplorer.com/football/england/championship/millwall-bristol-city/0Q90Liwg/"

# Row 66: was match #65 content, now gets match data previously stored at row 71
$ws.Cells.Item(66, 6).Value = "Preston"
$ws.Cells.Item(66, 7).Value = 2
$ws.Cells.Item(66, 8).Value = "Plymouth"
$ws.Cells.Item(66, 9).Value = 1
$ws.Cells.Item(66, 10).Value = 2
$ws.Cells.Item(66, 11).Value = "03/09/2023 16:12"
$ws.Cells.Item(66, 12).Value = 1.98
$ws.Cells.Item(66, 13).Value = "16/09/2023 15:48"
$ws.Cells.Item(66, 14).Value = 3.68
$ws.Cells.Item(66, 15).Value = "03/09/2023 16:12"
$ws.Cells.Item(66, 16).Value = 3.87
$ws.Cells.Item(66, 17).Value = "16/09/2023 15:59"
$ws.Cells.Item(66, 18).Value = 3.81
$ws.Cells.Item(66, 19).Value = "03/09/2023 16:12"
$ws.Cells.Item(66, 20).Value = 3.77
$ws.Cells.Item(66, 21).Value = "16/09/2023 15:59"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/england/championship/preston-plymouth/2FPAhd4r/"

# Row 67: was match #66 content, now gets match data previously stored at row 66
$ws.Cells.Item(67, 6).Value = "QPR"
$ws.Cells.Item(67, 7).Value = 1
$ws.Cells.Item(67, 8).Value = "Sunderland"
$ws.Cells.Item(67, 9).Value = 3
$ws.Cells.Item(67, 10).Value = 2.82
$ws.Cells.Item(67, 11).Value = "03/09/2023 15:42"
$ws.Cells.Item(67, 12).Value = 2.98
$ws.Cells.Item(67, 13).Value = "16/09/2023 15:47"
$ws.Cells.Item(67, 14).Value = 3.4
$ws.Cells.Item(67, 15).Value = "03/09/2023 15:42"
$ws.Cells.Item(67, 16).Value = 3.52
$ws.Cells.Item(67, 17).Value = "16/09/2023 15:47"
$ws.Cells.Item(67, 18).Value = 2.6
$ws.Cells.Item(67, 19).Value = "03/09/2023 15:42"
$ws.Cells.Item(67, 20).Value = 2.45
$ws.Cells.Item(67, 21).Value = "16/09/2023 15:47"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/england/championship/qpr-sunderland/00rEiGJl/"

# Row 68: was match #67 content, now gets match data previously stored at row 67
$ws.Cells.Item(68, 6).Value = "Huddersfield"
$ws.Cells.Item(68, 7).Value = 2
$ws.Cells.Item(68, 8).Value = "Rotherham"
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 1.97
$ws.Cells.Item(68, 11).Value = "03/09/2023 15:42"
$ws.Cells.Item(68, 12).Value = 1.86
$ws.Cells.Item(68, 13).Value = "16/09/2023 15:56"
$ws.Cells.Item(68, 14).Value = 3.58
$ws.Cells.Item(68, 15).Value = "03/09/2023 15:42"
$ws.Cells.Item(68, 16).Value = 3.68
$ws.Cells.Item(68, 17).Value = "16/09/2023 15:56"
$ws.Cells.Item(68, 18).Value = 4.04
$ws.Cells.Item(68, 19).Value = "03/09/2023 15:42"
$ws.Cells.Item(68, 20).Value = 4.52
$ws.Cells.Item(68, 21).Value = "16/09/2023 15:56"
$ws.Cells.Item(68, 22).Value = "https://www.betexplorer.com/football/england/championship/huddersfield-rotherham/AHpWavC8/"

# Row 71: was match #70 content, now gets match data previously stored at row 68
$ws.Cells.Item(71, 6).Value = "Norwich"
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = "Stoke"
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 1.78
$ws.Cells.Item(71, 11).Value = "03/09/2023 15:42"
$ws.Cells.Item(71, 12).Value = 1.93
$ws.Cells.Item(71, 13).Value = "16/09/2023 15:32"
$ws.Cells.Item(71, 14).Value = 3.95
$ws.Cells.Item(71, 15).Value = "03/09/2023 15:42"
$ws.Cells.Item(71, 16).Value = 3.75
$ws.Cells.Item(71, 17).Value = "16/09/2023 15:39"
$ws.Cells.Item(71, 18).Value = 4.5
$ws.Cells.Item(71, 19).Value = "03/09/2023 15:42"
$ws.Cells.Item(71, 20).Value = 4.09
$ws.Cells.Item(71, 21).Value = "16/09/2023 15:32"
$ws.Cells.Item(71, 22).Value = "https://www.betexplorer.com/football/england/championship/norwich-stoke-city/pzwqcxdR/"

# Row 74: was match #73 content, now gets match data previously stored at row 76
$ws.Cells.Item(74, 6).Value = "Preston"
$ws.Cells.Item(74, 7).Value = 2
$ws.Cells.Item(74, 8).Value = "Birmingham"
$ws.Cells.Item(74, 9).Value = 1
$ws.Cells.Item(74, 10).Value = 2.24
$ws.Cells.Item(74, 11).Value = "16/09/2023 17:13"
$ws.Cells.Item(74, 12).Value = 2.42
$ws.Cells.Item(74, 13).Value = "19/09/2023 20:41"
$ws.Cells.Item(74, 14).Value = 3.3
$ws.Cells.Item(74, 15).Value = "16/09/2023 17:13"
$ws.Cells.Item(74, 16).Value = 3.09
$ws.Cells.Item(74, 17).Value = "19/09/2023 20:41"
$ws.Cells.Item(74, 18).Value = 3.53
$ws.Cells.Item(74, 19).Value = "16/09/2023 17:13"
$ws.Cells.Item(74, 20).Value = 3.44
$ws.Cells.Item(74, 21).Value = "19/09/2023 20:41"
$ws.Cells.Item(74, 22).Value = "https://www.betexplorer.com/football/england/championship/preston-birmingham/YD0BdQMD/"

# Row 75: was match #74 content, now gets match data previously stored at row 74
$ws.Cells.Item(75, 6).Value = "Bristol City"
$ws.Cells.Item(75, 7).Value = 4
$ws.Cells.Item(75, 8).Value = "Plymouth"
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = 1.98
$ws.Cells.Item(75, 11).Value = "16/09/2023 17:13"
$ws.Cells.Item(75, 12).Value = 1.62
$ws.Cells.Item(75, 13).Value = "19/09/2023 20:40"
$ws.Cells.Item(75, 14).Value = 3.77
$ws.Cells.Item(75, 15).Value = "16/09/2023 17:13"
$ws.Cells.Item(75, 16).Value = 4.42
$ws.Cells.Item(75, 17).Value = "19/09/2023 20:41"
$ws.Cells.Item(75, 18).Value = 3.78
$ws.Cells.Item(75, 19).Value = "16/09/2023 17:13"
$ws.Cells.Item(75, 20).Value = 5.34
$ws.Cells.Item(75, 21).Value = "19/09/2023 20:41"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/england/championship/bristol-city-plymouth/UuzVmYID/"

# Row 76: was match #75 content, now gets match data previously stored at row 75
$ws.Cells.Item(76, 6).Value = "Cardiff"
$ws.Cells.Item(76, 7).Value = 3
$ws.Cells.Item(76, 8).Value = "Coventry"
$ws.Cells.Item(76, 9).Value = 2
$ws.Cells.Item(76, 10).Value = 2.67
$ws.Cells.Item(76, 11).Value = "16/09/2023 20:12"
$ws.Cells.Item(76, 12).Value = 2.77
$ws.Cells.Item(76, 13).Value = "19/09/2023 20:44"
$ws.Cells.Item(76, 14).Value = 3.37
$ws.Cells.Item(76, 15).Value = "16/09/2023 20:12"
$ws.Cells.Item(76, 16).Value = 3.39
$ws.Cells.Item(76, 17).Value = "19/09/2023 20:15"
$ws.Cells.Item(76, 18).Value = 2.76
$ws.Cells.Item(76, 19).Value = "16/09/2023 20:12"
$ws.Cells.Item(76, 20).Value = 2.69
$ws.Cells.Item(76, 21).Value = "19/09/2023 20:44"
$ws.Cells.Item(76, 22).Value = "https://www.betexplorer.com/football/england/championship/cardiff-coventry/neFvaY02/"

# ===================================================================
# Part 2: append 6 new match rows (80-85) discovered/added at the end
# Row 1 (header) style for Indice (col A) = s1 (bold, thin border, centered)
# Column E (data_partida) style = s2 (date-time display format)
# Both replicated here by copying format from an existing data row (row 2)
# rather than hand-rolling new style entries, to match the workbook's style table.
# ===================================================================

# Row 80 (Indice 79)
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(80, 1).PasteSpecial(-4122)
$ws.Cells.Item(80, 2).Value = "england"
$ws.Cells.Item(80, 3).Value = "championship"
$ws.Cells.Item(80, 4).Value = "2023-2024"
$ws.Cells.Item(80, 5).Value = 45189.86458333334
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(80, 5).PasteSpecial(-4122)
$ws.Cells.Item(80, 6).Value = "Millwall"
$ws.Cells.Item(80, 7).Value = 3
$ws.Cells.Item(80, 8).Value = "Rotherham"
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 1.76
$ws.Cells.Item(80, 11).Value = "17/09/2023 09:57"
$ws.Cells.Item(80, 12).Value = 1.67
$ws.Cells.Item(80, 13).Value = "20/09/2023 20:41"
$ws.Cells.Item(80, 14).Value = 3.79
$ws.Cells.Item(80, 15).Value = "17/09/2023 09:57"
$ws.Cells.Item(80, 16).Value = 3.91
$ws.Cells.Item(80, 17).Value = "20/09/2023 20:44"
$ws.Cells.Item(80, 18).Value = 4.87
$ws.Cells.Item(80, 19).Value = "17/09/2023 09:57"
$ws.Cells.Item(80, 20).Value = 5.65
$ws.Cells.Item(80, 21).Value = "20/09/2023 20:44"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/england/championship/millwall-rotherham/WSthzz6B/"

# Row 81 (Indice 80)
$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(81, 1).PasteSpecial(-4122)
$ws.Cells.Item(81, 2).Value = "england"
$ws.Cells.Item(81, 3).Value = "championship"
$ws.Cells.Item(81, 4).Value = "2023-2024"
$ws.Cells.Item(81, 5).Value = 45189.86458333334
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(81, 5).PasteSpecial(-4122)
$ws.Cells.Item(81, 6).Value = "Hull"
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = "Leeds"
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 3.46
$ws.Cells.Item(81, 11).Value = "16/09/2023 17:40"
$ws.Cells.Item(81, 12).Value = 3.83
$ws.Cells.Item(81, 13).Value = "20/09/2023 20:42"
$ws.Cells.Item(81, 14).Value = 3.69
$ws.Cells.Item(81, 15).Value = "16/09/2023 17:40"
$ws.Cells.Item(81, 16).Value = 3.8
$ws.Cells.Item(81, 17).Value = "20/09/2023 20:42"
$ws.Cells.Item(81, 18).Value = 2.11
$ws.Cells.Item(81, 19).Value = "16/09/2023 17:40"
$ws.Cells.Item(81, 20).Value = 1.98
$ws.Cells.Item(81, 21).Value = "20/09/2023 20:42"
$ws.Cells.Item(81, 22).Value = "https://www.betexplorer.com/football/england/championship/hull-city-leeds/x6slyGj5/"

# Row 82 (Indice 81)
$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(82, 1).PasteSpecial(-4122)
$ws.Cells.Item(82, 2).Value = "england"
$ws.Cells.Item(82, 3).Value = "championship"
$ws.Cells.Item(82, 4).Value = "2023-2024"
$ws.Cells.Item(82, 5).Value = 45189.86458333334
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(82, 5).PasteSpecial(-4122)
$ws.Cells.Item(82, 6).Value = "Watford"
$ws.Cells.Item(82, 7).Value = 2
$ws.Cells.Item(82, 8).Value = "West Brom"
$ws.Cells.Item(82, 9).Value = 2
$ws.Cells.Item(82, 10).Value = 2.39
$ws.Cells.Item(82, 11).Value = "16/09/2023 17:13"
$ws.Cells.Item(82, 12).Value = 2.27
$ws.Cells.Item(82, 13).Value = "20/09/2023 20:42"
$ws.Cells.Item(82, 14).Value = 3.49
$ws.Cells.Item(82, 15).Value = "16/09/2023 17:13"
$ws.Cells.Item(82, 16).Value = 3.38
$ws.Cells.Item(82, 17).Value = "20/09/2023 20:36"
$ws.Cells.Item(82, 18).Value = 3.05
$ws.Cells.Item(82, 19).Value = "16/09/2023 17:13"
$ws.Cells.Item(82, 20).Value = 3.43
$ws.Cells.Item(82, 21).Value = "20/09/2023 20:42"
$ws.Cells.Item(82, 22).Value = "https://www.betexplorer.com/football/england/championship/watford-west-brom/Sfz1YDzO/"

# Row 83 (Indice 82)
$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(83, 1).PasteSpecial(-4122)
$ws.Cells.Item(83, 2).Value = "england"
$ws.Cells.Item(83, 3).Value = "championship"
$ws.Cells.Item(83, 4).Value = "2023-2024"
$ws.Cells.Item(83, 5).Value = 45189.86458333334
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(83, 5).PasteSpecial(-4122)
$ws.Cells.Item(83, 6).Value = "Blackburn"
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = "Sunderland"
$ws.Cells.Item(83, 9).Value = 3
$ws.Cells.Item(83, 10).Value = 2.25
$ws.Cells.Item(83, 11).Value = "16/09/2023 17:13"
$ws.Cells.Item(83, 12).Value = 2.29
$ws.Cells.Item(83, 13).Value = "20/09/2023 20:41"
$ws.Cells.Item(83, 14).Value = 3.76
$ws.Cells.Item(83, 15).Value = "16/09/2023 17:13"
$ws.Cells.Item(83, 16).Value = 3.71
$ws.Cells.Item(83, 17).Value = "20/09/2023 20:36"
$ws.Cells.Item(83, 18).Value = 3.09
$ws.Cells.Item(83, 19).Value = "16/09/2023 17:13"
$ws.Cells.Item(83, 20).Value = 3.11
$ws.Cells.Item(83, 21).Value = "20/09/2023 20:41"
$ws.Cells.Item(83, 22).Value = "https://www.betexplorer.com/football/england/championship/blackburn-sunderland/I3RnwxMh/"

# Row 84 (Indice 83)
$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(84, 1).PasteSpecial(-4122)
$ws.Cells.Item(84, 2).Value = "england"
$ws.Cells.Item(84, 3).Value = "championship"
$ws.Cells.Item(84, 4).Value = "2023-2024"
$ws.Cells.Item(84, 5).Value = 45189.86458333334
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(84, 5).PasteSpecial(-4122)
$ws.Cells.Item(84, 6).Value = "Huddersfield"
$ws.Cells.Item(84, 7).Value = 2
$ws.Cells.Item(84, 8).Value = "Stoke"
$ws.Cells.Item(84, 9).Value = 2
$ws.Cells.Item(84, 10).Value = 2.84
$ws.Cells.Item(84, 11).Value = "16/09/2023 17:13"
$ws.Cells.Item(84, 12).Value = 3.13
$ws.Cells.Item(84, 13).Value = "20/09/2023 20:42"
$ws.Cells.Item(84, 14).Value = 3.33
$ws.Cells.Item(84, 15).Value = "16/09/2023 17:13"
$ws.Cells.Item(84, 16).Value = 3.38
$ws.Cells.Item(84, 17).Value = "20/09/2023 20:38"
$ws.Cells.Item(84, 18).Value = 2.66
$ws.Cells.Item(84, 19).Value = "16/09/2023 17:13"
$ws.Cells.Item(84, 20).Value = 2.42
$ws.Cells.Item(84, 21).Value = "20/09/2023 20:42"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/england/championship/huddersfield-stoke-city/bFrpxdya/"

# Row 85 (Indice 84)
$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(85, 1).PasteSpecial(-4122)
$ws.Cells.Item(85, 2).Value = "england"
$ws.Cells.Item(85, 3).Value = "championship"
$ws.Cells.Item(85, 4).Value = "2023-2024"
$ws.Cells.Item(85, 5).Value = 45189.875
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(85, 5).PasteSpecial(-4122)
$ws.Cells.Item(85, 6).Value = "Norwich"
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = "Leicester"
$ws.Cells.Item(85, 9).Value = 2
$ws.Cells.Item(85, 10).Value = 2.69
$ws.Cells.Item(85, 11).Value = "16/09/2023 17:13"
$ws.Cells.Item(85, 12).Value = 2.43
$ws.Cells.Item(85, 13).Value = "20/09/2023 20:58"
$ws.Cells.Item(85, 14).Value = 3.53
$ws.Cells.Item(85, 15).Value = "16/09/2023 17:13"
$ws.Cells.Item(85, 16).Value = 3.6
$ws.Cells.Item(85, 17).Value = "20/09/2023 20:58"
$ws.Cells.Item(85, 18).Value = 2.65
$ws.Cells.Item(85, 19).Value = "16/09/2023 17:13"
$ws.Cells.Item(85, 20).Value = 2.96
$ws.Cells.Item(85, 21).Value = "20/09/2023 20:58"
$ws.Cells.Item(85, 22).Value = "https://www.betexplorer.com/football/england/championship/norwich-leicester/pWpdZgLH/"

$excel.CutCopyMode = $false

